# examples updated for v17.01
# Populate A2:A7 with the new shared-string values and move the active
# selection to the last entry, matching the authored worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "pause"
$ws.Range("A3").Value = "structure"
$ws.Range("A4").Value = "paws"
$ws.Range("A5").Value = "pores"
$ws.Range("A6").Value = "pours"
$ws.Range("A7").Value = "braze"

$ws.Range("A7").Select() | Out-Null
